$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header row - C1 keeps the same text "The closest predicted product"
# (left untouched, value unchanged)

# Row 2: new reactant/product pair inserted here (pushing the former content down logically),
# and its "closest predicted product" cell (C2) is cleared (kept its formatting/style).
$ws.Range("A2").Value = "O=Cc1cccc(c1)OC(C)C"
$ws.Range("B2").Value = "OCc1cccc(c1)OC(C)C"
$ws.Range("C2").ClearContents()

# Row 3: new reactant/product pair; C3 had no special formatting, so clearing it
# removes the cell entirely.
$ws.Range("A3").Value = "O=Cc1c(C)cccc1C(C)C"
$ws.Range("A3").ClearFormats()
$ws.Range("B3").Value = "OCc1c(C)cccc1C(C)C"
$ws.Range("C3").ClearContents()

# Row 4: shifted-in values from what used to be row 2; C4 keeps formatting, cleared.
$ws.Range("A4").Value = "O=CCc1ccccc1"
$ws.Range("A4").ClearFormats()
$ws.Range("B4").Value = "OCCc1ccccc1"
$ws.Range("C4").ClearContents()

# Row 5: shifted-in values from what used to be row 3; C5 keeps formatting, cleared.
$ws.Range("A5").Value = "O=CCc1cccc(c1)C(=O)C"
$ws.Range("A5").ClearFormats()
$ws.Range("B5").Value = "OCCc1cccc(c1)C(=O)C"
$ws.Range("B5").ClearFormats()
$ws.Range("C5").ClearContents()

# Row 6: shifted-in values; C6 had no special formatting, so clearing removes the cell.
$ws.Range("A6").Value = "CCc1ccc(cc1)C(=O)C"
$ws.Range("A6").ClearFormats()
$ws.Range("B6").Value = "CCc1ccc(cc1)C(O)C"
$ws.Range("C6").ClearContents()

# Update the selected cell shown in the saved workbook view.
$ws.Range("B11").Select()

Write-Host "done"
